$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Leading apostrophes force Excel to keep purely-numeric-looking Price
# strings (e.g. "271.32") stored as text, matching the source data which
# never contains real numeric cells in this column.

$ws.Range("D2").Value = "43.697.26"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "2.311.20"
$ws.Range("E3").Value = "  +2.86%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'271.32"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("D6").Value = "'93.87"
$ws.Range("E6").Value = "  +7.17%  "
$ws.Range("E7").Value = "  +0.87%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.623"
$ws.Range("E9").Value = "  +1.07%  "
$ws.Range("D10").Value = "'44.93"
$ws.Range("E10").Value = "  -2.09%  "
$ws.Range("D11").Value = "'0.0940"
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").Value = "'8.18"
$ws.Range("E12").Value = "  +8.21%  "
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("D14").Value = "2.652.53"
$ws.Range("E14").Value = "  +2.64%  "
$ws.Range("D15").Value = "'15.35"
$ws.Range("E15").Value = "  +2.25%  "
$ws.Range("D16").Value = "'0.846"
$ws.Range("E16").Value = "  +5.79%  "
$ws.Range("D17").Value = "2.319.14"
$ws.Range("E17").Value = "  +3.09%  "
$ws.Range("D18").Value = "43.682.16"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("E19").Value = "  +1.65%  "
$ws.Range("D20").Value = "'6.28"
$ws.Range("E20").Value = "  +3.97%  "
$ws.Range("D21").Value = "'71.69"
$ws.Range("E21").Value = "  +1.95%  "
$ws.Range("D22").Value = "'239.26"
$ws.Range("E22").Value = "  +2.28%  "
$ws.Range("D23").Value = "'2.28"
$ws.Range("E23").Value = "  -5.25%  "
$ws.Range("E24").Value = "  +8.11%  "
$ws.Range("D26").Value = "'11.39"
$ws.Range("E26").Value = "  +3.36%  "
$ws.Range("D27").Value = "'2.51"
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("E28").Value = "  +4.84%  "
$ws.Range("D29").Value = "'3.39"
$ws.Range("E29").Value = "  -5.40%  "
$ws.Range("D30").Value = "'38.94"
$ws.Range("E30").Value = "  -4.29%  "
$ws.Range("D31").Value = "'22.62"
$ws.Range("E31").Value = "  +8.77%  "
$ws.Range("D32").Value = "'172.43"
$ws.Range("E32").Value = "  -1.77%  "
$ws.Range("D33").Value = "'0.0903"
$ws.Range("E33").Value = "  -1.20%  "
$ws.Range("E34").Value = "  +2.57%  "
$ws.Range("E35").Value = "  +1.88%  "
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("D37").Value = "'4.50"
$ws.Range("E37").Value = "  +2.30%  "
$ws.Range("D38").Value = "'0.0357"
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("D39").Value = "'3.46"
$ws.Range("E39").Value = "  +3.52%  "
$ws.Range("D40").Value = "'0.236"
$ws.Range("E40").Value = "  +14.97%  "
$ws.Range("E41").Value = "  +6.75%  "
$ws.Range("D42").Value = "'12.18"
$ws.Range("E42").Value = "  -4.64%  "
$ws.Range("E43").Value = "  +17.13%  "
$ws.Range("D44").Value = "'5.45"
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("D45").Value = "'61.79"
$ws.Range("E45").Value = "  -5.58%  "
$ws.Range("D46").Value = "'8.95"
$ws.Range("E46").Value = "  +6.75%  "
$ws.Range("E47").Value = "  +2.88%  "
$ws.Range("D48").Value = "'100.43"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("D50").Value = "2.529.42"
$ws.Range("E50").Value = "  +2.53%  "
$ws.Range("D51").Value = "'0.428"
$ws.Range("E51").Value = "  -3.29%  "
